$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set number format to Text for the data range to prevent Excel from
# auto-converting numeric-looking or date-looking strings
$ws.Range("A2:Q45").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 'Neusha Barakati, Rocio Zapata Bustos, Dawn K. Coletta, Paul Langlais, Lindsay N. Kohler, Moulun Luo, Janet L. Funk, Wayne T. Willis, Lawrence J. Mandarino'
$ws.Range("B2").Value = 'Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona, USA.; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona, USA.; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona, USA.; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona, USA.; Center for Disparities in Diabetes, Obesity, and Metabolism, University of Arizona, Health Sciences, Tucson, Arizona, USA.; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona, USA.; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona, USA.; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona, USA.; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona, USA.'
$ws.Range("C2").Value = 'https://openalex.org/W4307974510'
$ws.Range("D2").Value = 'Fuel Selection in Skeletal Muscle Exercising at Low Intensity; Reliance on Carbohydrate in Very Sedentary Individuals'
$ws.Range("E2").Value = '2023-02-01'
$ws.Range("F2").Value = 'Metabolic Syndrome and Related Disorders'
$ws.Range("G2").Value = 'Mary Ann Liebert, Inc.'
$ws.Range("H2").Value = 'https://doi.org/10.1089/met.2022.0062'
$ws.Range("I2").Value = 'N/A'
$ws.Range("J2").Value = 'N/A'
$ws.Range("K2").Value = 'closed'
$ws.Range("L2").Value = 'en'
$ws.Range("M2").Value = '3'
$ws.Range("N2").Value = '2023'
$ws.Range("O2").Value = 'https://pubmed.ncbi.nlm.nih.gov/36318809'
$ws.Range("P2").Value = 'https://doi.org/10.1089/met.2022.0062'
$ws.Range("Q2").Value = 'article'

# Row 3
$ws.Range("A3").Value = 'Moulun Luo, Yuanhong Xu, Jike Li, Dongxia Luo, Li Zhu, Yanxi Wu, Xiaodong Liu, Pengfei Wu'
$ws.Range("B3").Value = 'Infectious Disease Laboratory, Chengdu Public Health Clinical Center , Chengdu , 610061 , China; Clinical Laboratory, Chengdu Public Health Clinical Center , Chengdu , 610061 , China; Infectious Disease Laboratory, Chengdu Public Health Clinical Center , Chengdu , 610061 , China; Infectious Disease Laboratory, Chengdu Public Health Clinical Center , Chengdu , 610061 , China; Hepatology Clinic, Chengdu Public Health Clinical Center , Chengdu , 610061 , China; Infectious Disease Laboratory, Chengdu Public Health Clinical Center , Chengdu , 610061 , China; Clinical Laboratory, Chengdu Public Health Clinical Center , Chengdu , 610061 , China; Infectious Disease Laboratory, Chengdu Public Health Clinical Center , Chengdu , 610061 , China'
$ws.Range("C3").Value = 'https://openalex.org/W4378717005'
$ws.Range("D3").Value = 'Vitamin D protects intestines from liver cirrhosis-induced inflammation and oxidative stress by inhibiting the TLR4/MyD88/NF-κB signaling pathway'
$ws.Range("E3").Value = '2023-01-01'
$ws.Range("F3").Value = 'Open Medicine'
$ws.Range("G3").Value = 'De Gruyter Open'
$ws.Range("H3").Value = 'https://doi.org/10.1515/med-2023-0714'
$ws.Range("I3").Value = 'cc-by'
$ws.Range("J3").Value = 'publishedVersion'
$ws.Range("K3").Value = 'gold'
$ws.Range("L3").Value = 'en'
$ws.Range("M3").Value = '1'
$ws.Range("N3").Value = '2023'
$ws.Range("O3").Value = 'https://pubmed.ncbi.nlm.nih.gov/37273916'
$ws.Range("P3").Value = 'https://doi.org/10.1515/med-2023-0714'
$ws.Range("Q3").Value = 'article'

# Row 4
$ws.Range("A4").Value = 'Moulun Luo, Gaetano Santulli'
$ws.Range("B4").Value = 'Center for Disparities in Diabetes, Obesity and Metabolism, University of Arizona Health Sciences, Tucson, AZ, United States; Department of Molecular Pharmacology, Einstein-Mount Sinai Diabetes Research Center (ES-DRC), Einstein Institute for Aging Research, Institute for Neuroimmunology and Inflammation (INI), Albert Einstein College of Medicine, New York, NY, United States'
$ws.Range("C4").Value = 'https://openalex.org/W4381429717'
$ws.Range("D4").Value = 'Editorial: The link between obesity, type 2 diabetes, and mitochondria'
$ws.Range("E4").Value = '2023-06-20'
$ws.Range("F4").Value = 'Frontiers in Endocrinology'
$ws.Range("G4").Value = 'Frontiers Media'
$ws.Range("H4").Value = 'https://doi.org/10.3389/fendo.2023.1229935'
$ws.Range("I4").Value = 'cc-by'
$ws.Range("J4").Value = 'publishedVersion'
$ws.Range("K4").Value = 'gold'
$ws.Range("L4").Value = 'en'
$ws.Range("M4").Value = '1'
$ws.Range("N4").Value = '2023'
$ws.Range("O4").Value = 'https://pubmed.ncbi.nlm.nih.gov/37409237'
$ws.Range("P4").Value = 'https://doi.org/10.3389/fendo.2023.1229935'
$ws.Range("Q4").Value = 'editorial'

# Row 5
$ws.Range("A5").Value = 'Mandi J. Corenblum, Aiden McRobbie-Johnson, Emma Carruth, Kelsey Bernard, Moulun Luo, Lawrence J. Mandarino, Steven L. Peterson, Maria Sans-Fuentes, Dean Billheimer, Timothy Maley, Erika D. Eggers, Lalitha Madhavan'
$ws.Range("B5").Value = 'Department of Neurology, University of Arizona, Tucson, AZ; Physiological Sciences Graduate Program, University of Arizona, Tucson, AZ; Physiology Undergraduate Program, University of Arizona, Tucson, AZ; Physiological Sciences Graduate Program, University of Arizona, Tucson, AZ; Department of Medicine, University of Arizona, Tucson, AZ; Department of Medicine, University of Arizona, Tucson, AZ; Statistical Consulting Lab, BIO5 Institute, University of Arizona, Tucson, AZ; Statistical Consulting Lab, BIO5 Institute, University of Arizona, Tucson, AZ; Statistical Consulting Lab, BIO5 Institute, University of Arizona, Tucson, AZ; Physiological Sciences Graduate Program, University of Arizona, Tucson, AZ; Departments of Physiology and Biomedical Engineering, University of Arizona, Tucson, AZ; Department of Neurology, University of Arizona, Tucson, AZ; Evelyn F McKnight Brain Institute and BIO5 Institute, University of Arizona, Tucson, AZ'
$ws.Range("C5").Value = 'https://openalex.org/W4384023243'
$ws.Range("D5").Value = 'Parallel neurodegenerative phenotypes in sporadic Parkinson’s disease fibroblasts and midbrain dopamine neurons'
$ws.Range("E5").Value = '2023-10-01'
$ws.Range("F5").Value = 'Progress in Neurobiology'
$ws.Range("G5").Value = 'Elsevier BV'
$ws.Range("H5").Value = 'https://doi.org/10.1016/j.pneurobio.2023.102501'
$ws.Range("I5").Value = 'cc-by-nc-nd'
$ws.Range("J5").Value = 'publishedVersion'
$ws.Range("K5").Value = 'hybrid'
$ws.Range("L5").Value = 'en'
$ws.Range("M5").Value = '1'
$ws.Range("N5").Value = '2023'
$ws.Range("O5").Value = 'https://pubmed.ncbi.nlm.nih.gov/37451330'
$ws.Range("P5").Value = 'https://doi.org/10.1016/j.pneurobio.2023.102501'
$ws.Range("Q5").Value = 'article'

# Row 6
$ws.Range("A6").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B6").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C6").Value = 'https://openalex.org/W4362632760'
$ws.Range("D6").Value = 'Supplement 2. Next generation sequencing results for first round screen from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E6").Value = '2023-04-03'
$ws.Range("F6").Value = 'N/A'
$ws.Range("G6").Value = 'N/A'
$ws.Range("H6").Value = 'https://doi.org/10.1158/1535-7163.22509898.v1'
$ws.Range("I6").Value = 'cc-by'
$ws.Range("J6").Value = 'submittedVersion'
$ws.Range("K6").Value = 'bronze'
$ws.Range("L6").Value = 'en'
$ws.Range("M6").Value = '0'
$ws.Range("N6").Value = '2023'
$ws.Range("O6").Value = 'NA'
$ws.Range("P6").Value = 'https://doi.org/10.1158/1535-7163.22509898.v1'
$ws.Range("Q6").Value = 'article'

# Row 7
$ws.Range("A7").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B7").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C7").Value = 'https://openalex.org/W4362638052'
$ws.Range("D7").Value = 'Table S1: IC50 (nM) for M cells deleted PSMC subunits treated with BTZ from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E7").Value = '2023-04-03'
$ws.Range("F7").Value = 'N/A'
$ws.Range("G7").Value = 'N/A'
$ws.Range("H7").Value = 'https://doi.org/10.1158/1535-7163.22509883.v1'
$ws.Range("I7").Value = 'cc-by'
$ws.Range("J7").Value = 'submittedVersion'
$ws.Range("K7").Value = 'green'
$ws.Range("L7").Value = 'en'
$ws.Range("M7").Value = '0'
$ws.Range("N7").Value = '2023'
$ws.Range("O7").Value = 'NA'
$ws.Range("P7").Value = 'https://doi.org/10.1158/1535-7163.22509883.v1'
$ws.Range("Q7").Value = 'article'

# Row 8
$ws.Range("A8").Value = 'Mandi J. Corenblum, Aiden McRobbie-Johnson, Emma Carruth, Kelsey Bernard, Moulun Luo, Lawrence J. Mandarino, Sayeh Peterson, Dean Billheimer, Timothy Maley, Erika D. Eggers, Lalitha Madhavan'
$ws.Range("B8").Value = 'University of Arizona; University of Arizona; University of Arizona; University of Arizona; University of Arizona; University of Arizona; University of Arizona; University of Arizona; University of Arizona; University of Arizona; University of Arizona'
$ws.Range("C8").Value = 'https://openalex.org/W4320179215'
$ws.Range("D8").Value = 'Parallel Neurodegenerative Phenotypes in Sporadic Parkinson’s Disease Fibroblasts and Midbrain Dopamine Neurons'
$ws.Range("E8").Value = '2023-02-12'
$ws.Range("F8").Value = 'bioRxiv (Cold Spring Harbor Laboratory)'
$ws.Range("G8").Value = 'Cold Spring Harbor Laboratory'
$ws.Range("H8").Value = 'https://doi.org/10.1101/2023.02.10.527867'
$ws.Range("I8").Value = 'N/A'
$ws.Range("J8").Value = 'submittedVersion'
$ws.Range("K8").Value = 'green'
$ws.Range("L8").Value = 'en'
$ws.Range("M8").Value = '0'
$ws.Range("N8").Value = '2023'
$ws.Range("O8").Value = 'https://pubmed.ncbi.nlm.nih.gov/36798207'
$ws.Range("P8").Value = 'https://doi.org/10.1101/2023.02.10.527867'
$ws.Range("Q8").Value = 'article'

# Row 9
$ws.Range("A9").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B9").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C9").Value = 'https://openalex.org/W4361824079'
$ws.Range("D9").Value = 'Data from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E9").Value = '2023-03-31'
$ws.Range("F9").Value = 'N/A'
$ws.Range("G9").Value = 'N/A'
$ws.Range("H9").Value = 'https://doi.org/10.1158/0008-5472.c.6509109.v1'
$ws.Range("I9").Value = 'N/A'
$ws.Range("J9").Value = 'submittedVersion'
$ws.Range("K9").Value = 'closed'
$ws.Range("L9").Value = 'en'
$ws.Range("M9").Value = '0'
$ws.Range("N9").Value = '2023'
$ws.Range("O9").Value = 'NA'
$ws.Range("P9").Value = 'https://doi.org/10.1158/0008-5472.c.6509109.v1'
$ws.Range("Q9").Value = 'article'

# Row 10
$ws.Range("A10").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B10").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C10").Value = 'https://openalex.org/W4362346493'
$ws.Range("D10").Value = 'Supplementary table 2 from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E10").Value = '2023-03-31'
$ws.Range("F10").Value = 'N/A'
$ws.Range("G10").Value = 'N/A'
$ws.Range("H10").Value = 'https://doi.org/10.1158/0008-5472.22414911'
$ws.Range("I10").Value = 'cc-by'
$ws.Range("J10").Value = 'submittedVersion'
$ws.Range("K10").Value = 'bronze'
$ws.Range("L10").Value = 'en'
$ws.Range("M10").Value = '0'
$ws.Range("N10").Value = '2023'
$ws.Range("O10").Value = 'NA'
$ws.Range("P10").Value = 'https://doi.org/10.1158/0008-5472.22414911'
$ws.Range("Q10").Value = 'article'

# Row 11
$ws.Range("A11").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B11").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C11").Value = 'https://openalex.org/W4362348194'
$ws.Range("D11").Value = 'supplementary figures from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E11").Value = '2023-03-31'
$ws.Range("F11").Value = 'N/A'
$ws.Range("G11").Value = 'N/A'
$ws.Range("H11").Value = 'https://doi.org/10.1158/0008-5472.22414899'
$ws.Range("I11").Value = 'cc-by'
$ws.Range("J11").Value = 'submittedVersion'
$ws.Range("K11").Value = 'bronze'
$ws.Range("L11").Value = 'en'
$ws.Range("M11").Value = '0'
$ws.Range("N11").Value = '2023'
$ws.Range("O11").Value = 'NA'
$ws.Range("P11").Value = 'https://doi.org/10.1158/0008-5472.22414899'
$ws.Range("Q11").Value = 'article'

# Row 12
$ws.Range("A12").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B12").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C12").Value = 'https://openalex.org/W4362348196'
$ws.Range("D12").Value = 'Supplementary table 5 from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E12").Value = '2023-03-31'
$ws.Range("F12").Value = 'N/A'
$ws.Range("G12").Value = 'N/A'
$ws.Range("H12").Value = 'https://doi.org/10.1158/0008-5472.22414902'
$ws.Range("I12").Value = 'cc-by'
$ws.Range("J12").Value = 'submittedVersion'
$ws.Range("K12").Value = 'bronze'
$ws.Range("L12").Value = 'en'
$ws.Range("M12").Value = '0'
$ws.Range("N12").Value = '2023'
$ws.Range("O12").Value = 'NA'
$ws.Range("P12").Value = 'https://doi.org/10.1158/0008-5472.22414902'
$ws.Range("Q12").Value = 'article'

# Row 13
$ws.Range("A13").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B13").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C13").Value = 'https://openalex.org/W4362421757'
$ws.Range("D13").Value = 'Supplementary table 4 from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E13").Value = '2023-03-31'
$ws.Range("F13").Value = 'N/A'
$ws.Range("G13").Value = 'N/A'
$ws.Range("H13").Value = 'https://doi.org/10.1158/0008-5472.22414905.v1'
$ws.Range("I13").Value = 'cc-by'
$ws.Range("J13").Value = 'submittedVersion'
$ws.Range("K13").Value = 'green'
$ws.Range("L13").Value = 'en'
$ws.Range("M13").Value = '0'
$ws.Range("N13").Value = '2023'
$ws.Range("O13").Value = 'NA'
$ws.Range("P13").Value = 'https://doi.org/10.1158/0008-5472.22414905.v1'
$ws.Range("Q13").Value = 'article'

# Row 14
$ws.Range("A14").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B14").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C14").Value = 'https://openalex.org/W4362421849'
$ws.Range("D14").Value = 'supplementary figures from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E14").Value = '2023-03-31'
$ws.Range("F14").Value = 'N/A'
$ws.Range("G14").Value = 'N/A'
$ws.Range("H14").Value = 'https://doi.org/10.1158/0008-5472.22414899.v1'
$ws.Range("I14").Value = 'cc-by'
$ws.Range("J14").Value = 'submittedVersion'
$ws.Range("K14").Value = 'green'
$ws.Range("L14").Value = 'en'
$ws.Range("M14").Value = '0'
$ws.Range("N14").Value = '2023'
$ws.Range("O14").Value = 'NA'
$ws.Range("P14").Value = 'https://doi.org/10.1158/0008-5472.22414899.v1'
$ws.Range("Q14").Value = 'article'

# Row 15
$ws.Range("A15").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B15").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C15").Value = 'https://openalex.org/W4362421881'
$ws.Range("D15").Value = 'Spplementary table 1 from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E15").Value = '2023-03-31'
$ws.Range("F15").Value = 'N/A'
$ws.Range("G15").Value = 'N/A'
$ws.Range("H15").Value = 'https://doi.org/10.1158/0008-5472.22414914.v1'
$ws.Range("I15").Value = 'cc-by'
$ws.Range("J15").Value = 'submittedVersion'
$ws.Range("K15").Value = 'green'
$ws.Range("L15").Value = 'en'
$ws.Range("M15").Value = '0'
$ws.Range("N15").Value = '2023'
$ws.Range("O15").Value = 'NA'
$ws.Range("P15").Value = 'https://doi.org/10.1158/0008-5472.22414914.v1'
$ws.Range("Q15").Value = 'article'

# Row 16
$ws.Range("A16").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B16").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C16").Value = 'https://openalex.org/W4362422092'
$ws.Range("D16").Value = 'Supplementary table 3 from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E16").Value = '2023-03-31'
$ws.Range("F16").Value = 'N/A'
$ws.Range("G16").Value = 'N/A'
$ws.Range("H16").Value = 'https://doi.org/10.1158/0008-5472.22414908.v1'
$ws.Range("I16").Value = 'cc-by'
$ws.Range("J16").Value = 'submittedVersion'
$ws.Range("K16").Value = 'green'
$ws.Range("L16").Value = 'en'
$ws.Range("M16").Value = '0'
$ws.Range("N16").Value = '2023'
$ws.Range("O16").Value = 'NA'
$ws.Range("P16").Value = 'https://doi.org/10.1158/0008-5472.22414908.v1'
$ws.Range("Q16").Value = 'article'

# Row 17
$ws.Range("A17").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B17").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C17").Value = 'https://openalex.org/W4362422103'
$ws.Range("D17").Value = 'Supplementary table 2 from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E17").Value = '2023-03-31'
$ws.Range("F17").Value = 'N/A'
$ws.Range("G17").Value = 'N/A'
$ws.Range("H17").Value = 'https://doi.org/10.1158/0008-5472.22414911.v1'
$ws.Range("I17").Value = 'cc-by'
$ws.Range("J17").Value = 'submittedVersion'
$ws.Range("K17").Value = 'green'
$ws.Range("L17").Value = 'en'
$ws.Range("M17").Value = '0'
$ws.Range("N17").Value = '2023'
$ws.Range("O17").Value = 'NA'
$ws.Range("P17").Value = 'https://doi.org/10.1158/0008-5472.22414911.v1'
$ws.Range("Q17").Value = 'article'

# Row 18
$ws.Range("A18").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B18").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C18").Value = 'https://openalex.org/W4362422259'
$ws.Range("D18").Value = 'Supplementary table 5 from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E18").Value = '2023-03-31'
$ws.Range("F18").Value = 'N/A'
$ws.Range("G18").Value = 'N/A'
$ws.Range("H18").Value = 'https://doi.org/10.1158/0008-5472.22414902.v1'
$ws.Range("I18").Value = 'cc-by'
$ws.Range("J18").Value = 'submittedVersion'
$ws.Range("K18").Value = 'green'
$ws.Range("L18").Value = 'en'
$ws.Range("M18").Value = '0'
$ws.Range("N18").Value = '2023'
$ws.Range("O18").Value = 'NA'
$ws.Range("P18").Value = 'https://doi.org/10.1158/0008-5472.22414902.v1'
$ws.Range("Q18").Value = 'article'

# Row 19
$ws.Range("A19").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B19").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C19").Value = 'https://openalex.org/W4362446471'
$ws.Range("D19").Value = 'Supplementary table 4 from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E19").Value = '2023-03-31'
$ws.Range("F19").Value = 'N/A'
$ws.Range("G19").Value = 'N/A'
$ws.Range("H19").Value = 'https://doi.org/10.1158/0008-5472.22414905'
$ws.Range("I19").Value = 'cc-by'
$ws.Range("J19").Value = 'submittedVersion'
$ws.Range("K19").Value = 'bronze'
$ws.Range("L19").Value = 'en'
$ws.Range("M19").Value = '0'
$ws.Range("N19").Value = '2023'
$ws.Range("O19").Value = 'NA'
$ws.Range("P19").Value = 'https://doi.org/10.1158/0008-5472.22414905'
$ws.Range("Q19").Value = 'article'

# Row 20
$ws.Range("A20").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B20").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C20").Value = 'https://openalex.org/W4362461672'
$ws.Range("D20").Value = 'Supplementary table 3 from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E20").Value = '2023-03-31'
$ws.Range("F20").Value = 'N/A'
$ws.Range("G20").Value = 'N/A'
$ws.Range("H20").Value = 'https://doi.org/10.1158/0008-5472.22414908'
$ws.Range("I20").Value = 'cc-by'
$ws.Range("J20").Value = 'submittedVersion'
$ws.Range("K20").Value = 'bronze'
$ws.Range("L20").Value = 'en'
$ws.Range("M20").Value = '0'
$ws.Range("N20").Value = '2023'
$ws.Range("O20").Value = 'NA'
$ws.Range("P20").Value = 'https://doi.org/10.1158/0008-5472.22414908'
$ws.Range("Q20").Value = 'article'

# Row 21
$ws.Range("A21").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B21").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C21").Value = 'https://openalex.org/W4362461723'
$ws.Range("D21").Value = 'Spplementary table 1 from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E21").Value = '2023-03-31'
$ws.Range("F21").Value = 'N/A'
$ws.Range("G21").Value = 'N/A'
$ws.Range("H21").Value = 'https://doi.org/10.1158/0008-5472.22414914'
$ws.Range("I21").Value = 'cc-by'
$ws.Range("J21").Value = 'submittedVersion'
$ws.Range("K21").Value = 'bronze'
$ws.Range("L21").Value = 'en'
$ws.Range("M21").Value = '0'
$ws.Range("N21").Value = '2023'
$ws.Range("O21").Value = 'NA'
$ws.Range("P21").Value = 'https://doi.org/10.1158/0008-5472.22414914'
$ws.Range("Q21").Value = 'article'

# Row 22
$ws.Range("A22").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B22").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C22").Value = 'https://openalex.org/W4362528337'
$ws.Range("D22").Value = 'Table S1: IC50 (nM) for M cells deleted PSMC subunits treated with BTZ from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E22").Value = '2023-04-03'
$ws.Range("F22").Value = 'N/A'
$ws.Range("G22").Value = 'N/A'
$ws.Range("H22").Value = 'https://doi.org/10.1158/1535-7163.22509883'
$ws.Range("I22").Value = 'cc-by'
$ws.Range("J22").Value = 'submittedVersion'
$ws.Range("K22").Value = 'green'
$ws.Range("L22").Value = 'en'
$ws.Range("M22").Value = '0'
$ws.Range("N22").Value = '2023'
$ws.Range("O22").Value = 'NA'
$ws.Range("P22").Value = 'https://doi.org/10.1158/1535-7163.22509883'
$ws.Range("Q22").Value = 'article'

# Row 23
$ws.Range("A23").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B23").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C23").Value = 'https://openalex.org/W4362544870'
$ws.Range("D23").Value = 'Data from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E23").Value = '2023-04-03'
$ws.Range("F23").Value = 'N/A'
$ws.Range("G23").Value = 'N/A'
$ws.Range("H23").Value = 'https://doi.org/10.1158/1535-7163.c.6539428'
$ws.Range("I23").Value = 'N/A'
$ws.Range("J23").Value = 'submittedVersion'
$ws.Range("K23").Value = 'green'
$ws.Range("L23").Value = 'en'
$ws.Range("M23").Value = '0'
$ws.Range("N23").Value = '2023'
$ws.Range("O23").Value = 'NA'
$ws.Range("P23").Value = 'https://doi.org/10.1158/1535-7163.c.6539428'
$ws.Range("Q23").Value = 'article'

# Row 24
$ws.Range("A24").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B24").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C24").Value = 'https://openalex.org/W4362548392'
$ws.Range("D24").Value = 'FigureS2 from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E24").Value = '2023-04-03'
$ws.Range("F24").Value = 'N/A'
$ws.Range("G24").Value = 'N/A'
$ws.Range("H24").Value = 'https://doi.org/10.1158/1535-7163.22509904'
$ws.Range("I24").Value = 'cc-by'
$ws.Range("J24").Value = 'submittedVersion'
$ws.Range("K24").Value = 'green'
$ws.Range("L24").Value = 'en'
$ws.Range("M24").Value = '0'
$ws.Range("N24").Value = '2023'
$ws.Range("O24").Value = 'NA'
$ws.Range("P24").Value = 'https://doi.org/10.1158/1535-7163.22509904'
$ws.Range("Q24").Value = 'article'

# Row 25
$ws.Range("A25").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B25").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C25").Value = 'https://openalex.org/W4362548404'
$ws.Range("D25").Value = 'Supplemental Figure legends from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E25").Value = '2023-04-03'
$ws.Range("F25").Value = 'N/A'
$ws.Range("G25").Value = 'N/A'
$ws.Range("H25").Value = 'https://doi.org/10.1158/1535-7163.22509886'
$ws.Range("I25").Value = 'cc-by'
$ws.Range("J25").Value = 'submittedVersion'
$ws.Range("K25").Value = 'green'
$ws.Range("L25").Value = 'da'
$ws.Range("M25").Value = '0'
$ws.Range("N25").Value = '2023'
$ws.Range("O25").Value = 'NA'
$ws.Range("P25").Value = 'https://doi.org/10.1158/1535-7163.22509886'
$ws.Range("Q25").Value = 'article'

# Row 26
$ws.Range("A26").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B26").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C26").Value = 'https://openalex.org/W4362548466'
$ws.Range("D26").Value = 'Supplement 2. Next generation sequencing results for first round screen from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E26").Value = '2023-04-03'
$ws.Range("F26").Value = 'N/A'
$ws.Range("G26").Value = 'N/A'
$ws.Range("H26").Value = 'https://doi.org/10.1158/1535-7163.22509898'
$ws.Range("I26").Value = 'cc-by'
$ws.Range("J26").Value = 'submittedVersion'
$ws.Range("K26").Value = 'green'
$ws.Range("L26").Value = 'en'
$ws.Range("M26").Value = '0'
$ws.Range("N26").Value = '2023'
$ws.Range("O26").Value = 'NA'
$ws.Range("P26").Value = 'https://doi.org/10.1158/1535-7163.22509898'
$ws.Range("Q26").Value = 'article'

# Row 27
$ws.Range("A27").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B27").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C27").Value = 'https://openalex.org/W4362548474'
$ws.Range("D27").Value = 'Figure S1 from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E27").Value = '2023-04-03'
$ws.Range("F27").Value = 'N/A'
$ws.Range("G27").Value = 'N/A'
$ws.Range("H27").Value = 'https://doi.org/10.1158/1535-7163.22509907'
$ws.Range("I27").Value = 'cc-by'
$ws.Range("J27").Value = 'submittedVersion'
$ws.Range("K27").Value = 'green'
$ws.Range("L27").Value = 'en'
$ws.Range("M27").Value = '0'
$ws.Range("N27").Value = '2023'
$ws.Range("O27").Value = 'NA'
$ws.Range("P27").Value = 'https://doi.org/10.1158/1535-7163.22509907'
$ws.Range("Q27").Value = 'article'

# Row 28
$ws.Range("A28").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B28").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C28").Value = 'https://openalex.org/W4362548485'
$ws.Range("D28").Value = 'Supplement 3. sgRNA sequences used for second round screen. from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E28").Value = '2023-04-03'
$ws.Range("F28").Value = 'N/A'
$ws.Range("G28").Value = 'N/A'
$ws.Range("H28").Value = 'https://doi.org/10.1158/1535-7163.22509895'
$ws.Range("I28").Value = 'cc-by'
$ws.Range("J28").Value = 'submittedVersion'
$ws.Range("K28").Value = 'green'
$ws.Range("L28").Value = 'en'
$ws.Range("M28").Value = '0'
$ws.Range("N28").Value = '2023'
$ws.Range("O28").Value = 'NA'
$ws.Range("P28").Value = 'https://doi.org/10.1158/1535-7163.22509895'
$ws.Range("Q28").Value = 'article'

# Row 29
$ws.Range("A29").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B29").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C29").Value = 'https://openalex.org/W4362548494'
$ws.Range("D29").Value = 'Supplement 1. sgRNA sequences targeting PSMC1 to C6 from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E29").Value = '2023-04-03'
$ws.Range("F29").Value = 'N/A'
$ws.Range("G29").Value = 'N/A'
$ws.Range("H29").Value = 'https://doi.org/10.1158/1535-7163.22509901'
$ws.Range("I29").Value = 'N/A'
$ws.Range("J29").Value = 'submittedVersion'
$ws.Range("K29").Value = 'closed'
$ws.Range("L29").Value = 'en'
$ws.Range("M29").Value = '0'
$ws.Range("N29").Value = '2023'
$ws.Range("O29").Value = 'NA'
$ws.Range("P29").Value = 'https://doi.org/10.1158/1535-7163.22509901'
$ws.Range("Q29").Value = 'article'

# Row 30
$ws.Range("A30").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B30").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C30").Value = 'https://openalex.org/W4362548496'
$ws.Range("D30").Value = 'Supplement 5. A total of 36 mutations in 19S proteasome subunits out of 895 patients were identified in CoMMpass study. from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E30").Value = '2023-04-03'
$ws.Range("F30").Value = 'N/A'
$ws.Range("G30").Value = 'N/A'
$ws.Range("H30").Value = 'https://doi.org/10.1158/1535-7163.22509889'
$ws.Range("I30").Value = 'cc-by'
$ws.Range("J30").Value = 'submittedVersion'
$ws.Range("K30").Value = 'green'
$ws.Range("L30").Value = 'en'
$ws.Range("M30").Value = '0'
$ws.Range("N30").Value = '2023'
$ws.Range("O30").Value = 'NA'
$ws.Range("P30").Value = 'https://doi.org/10.1158/1535-7163.22509889'
$ws.Range("Q30").Value = 'article'

# Row 31
$ws.Range("A31").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B31").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C31").Value = 'https://openalex.org/W4362548508'
$ws.Range("D31").Value = 'Supplement 4. Next generation sequencing results for second round screen. from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E31").Value = '2023-04-03'
$ws.Range("F31").Value = 'N/A'
$ws.Range("G31").Value = 'N/A'
$ws.Range("H31").Value = 'https://doi.org/10.1158/1535-7163.22509892'
$ws.Range("I31").Value = 'N/A'
$ws.Range("J31").Value = 'submittedVersion'
$ws.Range("K31").Value = 'green'
$ws.Range("L31").Value = 'en'
$ws.Range("M31").Value = '0'
$ws.Range("N31").Value = '2023'
$ws.Range("O31").Value = 'NA'
$ws.Range("P31").Value = 'https://doi.org/10.1158/1535-7163.22509892'
$ws.Range("Q31").Value = 'article'

# Row 32
$ws.Range("A32").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B32").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C32").Value = 'https://openalex.org/W4362568194'
$ws.Range("D32").Value = 'Data from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E32").Value = '2023-04-03'
$ws.Range("F32").Value = 'N/A'
$ws.Range("G32").Value = 'N/A'
$ws.Range("H32").Value = 'https://doi.org/10.1158/1535-7163.c.6539428.v1'
$ws.Range("I32").Value = 'N/A'
$ws.Range("J32").Value = 'submittedVersion'
$ws.Range("K32").Value = 'green'
$ws.Range("L32").Value = 'en'
$ws.Range("M32").Value = '0'
$ws.Range("N32").Value = '2023'
$ws.Range("O32").Value = 'NA'
$ws.Range("P32").Value = 'https://doi.org/10.1158/1535-7163.c.6539428.v1'
$ws.Range("Q32").Value = 'article'

# Row 33
$ws.Range("A33").Value = 'Hao Chen, Moulun Luo, Xiangping Wang, Ting Liang, Chung-Chi Huang, Chien-Hung Huang, Lining Wei'
$ws.Range("B33").Value = 'Department of Oncology, The Second Nanning People’s Hospital, Nanning, China; Department of Oncology, The Second Nanning People’s Hospital, Nanning, China; Department of Oncology, The Second Nanning People’s Hospital, Nanning, China; Department of Oncology, The Second Nanning People’s Hospital, Nanning, China; Department of Oncology, The Second Nanning People’s Hospital, Nanning, China; Department of Oncology, The Second Nanning People’s Hospital, No. 13 Dancun Road, Jiangnan District, Nanning, 530031, Guangxi, China; Department of Endoscopy, The Affiliated Tumor Hospital of Guangxi Medical University, Nanning, China'
$ws.Range("C33").Value = 'https://openalex.org/W4362601969'
$ws.Range("D33").Value = 'Correction: Inhibition of PAD4 enhances radiosensitivity and inhibits aggressive phenotypes of nasopharyngeal carcinoma cells'
$ws.Range("E33").Value = '2023-04-04'
$ws.Range("F33").Value = 'Cellular & Molecular Biology Letters'
$ws.Range("G33").Value = 'BioMed Central'
$ws.Range("H33").Value = 'https://doi.org/10.1186/s11658-023-00444-x'
$ws.Range("I33").Value = 'cc-by'
$ws.Range("J33").Value = 'publishedVersion'
$ws.Range("K33").Value = 'gold'
$ws.Range("L33").Value = 'en'
$ws.Range("M33").Value = '0'
$ws.Range("N33").Value = '2023'
$ws.Range("O33").Value = 'https://pubmed.ncbi.nlm.nih.gov/37016289'
$ws.Range("P33").Value = 'https://doi.org/10.1186/s11658-023-00444-x'
$ws.Range("Q33").Value = 'article'

# Row 34
$ws.Range("A34").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B34").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C34").Value = 'https://openalex.org/W4362633184'
$ws.Range("D34").Value = 'Supplemental Figure legends from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E34").Value = '2023-04-03'
$ws.Range("F34").Value = 'N/A'
$ws.Range("G34").Value = 'N/A'
$ws.Range("H34").Value = 'https://doi.org/10.1158/1535-7163.22509886.v1'
$ws.Range("I34").Value = 'cc-by'
$ws.Range("J34").Value = 'submittedVersion'
$ws.Range("K34").Value = 'green'
$ws.Range("L34").Value = 'da'
$ws.Range("M34").Value = '0'
$ws.Range("N34").Value = '2023'
$ws.Range("O34").Value = 'NA'
$ws.Range("P34").Value = 'https://doi.org/10.1158/1535-7163.22509886.v1'
$ws.Range("Q34").Value = 'article'

# Row 35
$ws.Range("A35").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B35").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C35").Value = 'https://openalex.org/W4362634321'
$ws.Range("D35").Value = 'Figure S1 from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E35").Value = '2023-04-03'
$ws.Range("F35").Value = 'N/A'
$ws.Range("G35").Value = 'N/A'
$ws.Range("H35").Value = 'https://doi.org/10.1158/1535-7163.22509907.v1'
$ws.Range("I35").Value = 'cc-by'
$ws.Range("J35").Value = 'submittedVersion'
$ws.Range("K35").Value = 'green'
$ws.Range("L35").Value = 'en'
$ws.Range("M35").Value = '0'
$ws.Range("N35").Value = '2023'
$ws.Range("O35").Value = 'NA'
$ws.Range("P35").Value = 'https://doi.org/10.1158/1535-7163.22509907.v1'
$ws.Range("Q35").Value = 'article'

# Row 36
$ws.Range("A36").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B36").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C36").Value = 'https://openalex.org/W4362635423'
$ws.Range("D36").Value = 'FigureS2 from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E36").Value = '2023-04-03'
$ws.Range("F36").Value = 'N/A'
$ws.Range("G36").Value = 'N/A'
$ws.Range("H36").Value = 'https://doi.org/10.1158/1535-7163.22509904.v1'
$ws.Range("I36").Value = 'cc-by'
$ws.Range("J36").Value = 'submittedVersion'
$ws.Range("K36").Value = 'green'
$ws.Range("L36").Value = 'en'
$ws.Range("M36").Value = '0'
$ws.Range("N36").Value = '2023'
$ws.Range("O36").Value = 'NA'
$ws.Range("P36").Value = 'https://doi.org/10.1158/1535-7163.22509904.v1'
$ws.Range("Q36").Value = 'article'

# Row 37
$ws.Range("A37").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B37").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C37").Value = 'https://openalex.org/W4362636285'
$ws.Range("D37").Value = 'Supplement 5. A total of 36 mutations in 19S proteasome subunits out of 895 patients were identified in CoMMpass study. from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E37").Value = '2023-04-03'
$ws.Range("F37").Value = 'N/A'
$ws.Range("G37").Value = 'N/A'
$ws.Range("H37").Value = 'https://doi.org/10.1158/1535-7163.22509889.v1'
$ws.Range("I37").Value = 'cc-by'
$ws.Range("J37").Value = 'submittedVersion'
$ws.Range("K37").Value = 'green'
$ws.Range("L37").Value = 'en'
$ws.Range("M37").Value = '0'
$ws.Range("N37").Value = '2023'
$ws.Range("O37").Value = 'NA'
$ws.Range("P37").Value = 'https://doi.org/10.1158/1535-7163.22509889.v1'
$ws.Range("Q37").Value = 'article'

# Row 38
$ws.Range("A38").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B38").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C38").Value = 'https://openalex.org/W4362637849'
$ws.Range("D38").Value = 'Supplement 1. sgRNA sequences targeting PSMC1 to C6 from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E38").Value = '2023-04-03'
$ws.Range("F38").Value = 'N/A'
$ws.Range("G38").Value = 'N/A'
$ws.Range("H38").Value = 'https://doi.org/10.1158/1535-7163.22509901.v1'
$ws.Range("I38").Value = 'cc-by'
$ws.Range("J38").Value = 'submittedVersion'
$ws.Range("K38").Value = 'green'
$ws.Range("L38").Value = 'en'
$ws.Range("M38").Value = '0'
$ws.Range("N38").Value = '2023'
$ws.Range("O38").Value = 'NA'
$ws.Range("P38").Value = 'https://doi.org/10.1158/1535-7163.22509901.v1'
$ws.Range("Q38").Value = 'article'

# Row 39
$ws.Range("A39").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert A. Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B39").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C39").Value = 'https://openalex.org/W4362637961'
$ws.Range("D39").Value = 'Supplement 3. sgRNA sequences used for second round screen. from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E39").Value = '2023-04-03'
$ws.Range("F39").Value = 'N/A'
$ws.Range("G39").Value = 'N/A'
$ws.Range("H39").Value = 'https://doi.org/10.1158/1535-7163.22509895.v1'
$ws.Range("I39").Value = 'cc-by'
$ws.Range("J39").Value = 'submittedVersion'
$ws.Range("K39").Value = 'green'
$ws.Range("L39").Value = 'en'
$ws.Range("M39").Value = '0'
$ws.Range("N39").Value = '2023'
$ws.Range("O39").Value = 'NA'
$ws.Range("P39").Value = 'https://doi.org/10.1158/1535-7163.22509895.v1'
$ws.Range("Q39").Value = 'article'

# Row 40
$ws.Range("A40").Value = 'Chang‐Xin Shi, K. Martin Kortüm, Yuan Xiao Zhu, Laura A. Bruins, Patrick Jedlowski, Patrick G. Votruba, Moulun Luo, Robert Stewart, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B40").Value = '; ; ; ; ; ; ; ; ; ; '
$ws.Range("C40").Value = 'https://openalex.org/W4362639398'
$ws.Range("D40").Value = 'Supplement 4. Next generation sequencing results for second round screen. from CRISPR Genome-Wide Screening Identifies Dependence on the Proteasome Subunit PSMC6 for Bortezomib Sensitivity in Multiple Myeloma'
$ws.Range("E40").Value = '2023-04-03'
$ws.Range("F40").Value = 'N/A'
$ws.Range("G40").Value = 'N/A'
$ws.Range("H40").Value = 'https://doi.org/10.1158/1535-7163.22509892.v1'
$ws.Range("I40").Value = 'cc-by'
$ws.Range("J40").Value = 'submittedVersion'
$ws.Range("K40").Value = 'green'
$ws.Range("L40").Value = 'en'
$ws.Range("M40").Value = '0'
$ws.Range("N40").Value = '2023'
$ws.Range("O40").Value = 'NA'
$ws.Range("P40").Value = 'https://doi.org/10.1158/1535-7163.22509892.v1'
$ws.Range("Q40").Value = 'article'

# Row 41
$ws.Range("A41").Value = 'Yuanmin Wang, Moulun Luo, Fukun Zhao, Mu Chun Su, Meiyan Li'
$ws.Range("B41").Value = 'Department of pharmacy, The First People''s Hospital of Zunyi (The Third Affiliated Hospital of Zunyi Medical University), China; Department of pharmacy, The First People''s Hospital of Zunyi (The Third Affiliated Hospital of Zunyi Medical University), China; Department of pharmacy, The First People''s Hospital of Zunyi (The Third Affiliated Hospital of Zunyi Medical University), China; Department of pharmacy, The First People''s Hospital of Zunyi (The Third Affiliated Hospital of Zunyi Medical University), China; Department of pharmacy, The First People''s Hospital of Zunyi (The Third Affiliated Hospital of Zunyi Medical University), China'
$ws.Range("C41").Value = 'https://openalex.org/W4378804477'
$ws.Range("D41").Value = 'Bioinformatics analysis of key genes in patients with sarcoidosis and prediction of traditional Chinese Medicine'
$ws.Range("E41").Value = '2023-02-10'
$ws.Range("F41").Value = 'N/A'
$ws.Range("G41").Value = 'N/A'
$ws.Range("H41").Value = 'https://doi.org/10.1145/3592686.3592739'
$ws.Range("I41").Value = 'N/A'
$ws.Range("J41").Value = 'N/A'
$ws.Range("K41").Value = 'closed'
$ws.Range("L41").Value = 'en'
$ws.Range("M41").Value = '0'
$ws.Range("N41").Value = '2023'
$ws.Range("O41").Value = 'NA'
$ws.Range("P41").Value = 'https://doi.org/10.1145/3592686.3592739'
$ws.Range("Q41").Value = 'article'

# Row 42
$ws.Range("A42").Value = 'Yi Liu, X. Ai, G. C. Xiao, Yaxuan Li, Wu Ling-Hui, Liangliang Wang, J. Dong, M. Y. Dong, Q. X. Geng, Moulun Luo, Na Yan, Andrew Wang, Chenxu Wang, Meng Wang, Lei Zhang, Liang Zhang, Ruikai Zhang, Yao Zhang, M. G. Zhao, Yang Zhou'
$ws.Range("B42").Value = '; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; '
$ws.Range("C42").Value = 'https://openalex.org/W4386557787'
$ws.Range("D42").Value = 'Simulation study of BESIII with stitched CMOS pixel detector using ACTS'
$ws.Range("E42").Value = '2023-09-05'
$ws.Range("F42").Value = 'arXiv (Cornell University)'
$ws.Range("G42").Value = 'Cornell University'
$ws.Range("H42").Value = 'https://arxiv.org/abs/2309.02358'
$ws.Range("I42").Value = 'N/A'
$ws.Range("J42").Value = 'submittedVersion'
$ws.Range("K42").Value = 'green'
$ws.Range("L42").Value = 'en'
$ws.Range("M42").Value = '0'
$ws.Range("N42").Value = '2023'
$ws.Range("O42").Value = 'NA'
$ws.Range("P42").Value = 'https://doi.org/10.48550/arxiv.2309.02358'
$ws.Range("Q42").Value = 'article'

# Row 43
$ws.Range("A43").Value = 'Yi Liu, X. Ai, G. Y. Xiao, Yaxuan Li, Wu Ling-Hui, Liang-Liang Wang, J. Dong, Ming-Yi Dong, Q. X. Geng, Moulun Luo, Na Yan, Andrew Wang, Chenxu Wang, Meng Wang, Lei Zhang, Liang Zhang, Rui-Kai Zhang, Yao Zhang, M. G. Zhao, Yang Zhou'
$ws.Range("B43").Value = 'School of Physics and Microelectronics, Zhengzhou University, Zhengzhou, 450001, China; School of Physics and Microelectronics, Zhengzhou University, Zhengzhou, 450001, China; School of Physics, Nanjing University, Nanjing, 210093, China; School of Physics, Nankai University, Tianjin, 300071, China; Institute of High Energy Physics, Chinese Academy of Sciences, 19B Yuquan Road, Shijingshan District, Beijing, 100049, China; Institute of High Energy Physics, Chinese Academy of Sciences, 19B Yuquan Road, Shijingshan District, Beijing, 100049, China; Research Center for Particle Science and Technology, Institute of Frontier and Interdisciplinary Science, Shandong University, Qingdao, 266237, China; Institute of High Energy Physics, Chinese Academy of Sciences, 19B Yuquan Road, Shijingshan District, Beijing, 100049, China; Research Center for Particle Science and Technology, Institute of Frontier and Interdisciplinary Science, Shandong University, Qingdao, 266237, China; School of Information Science and Engineering, Harbin Institute of Technology, Weihai, 264209, China; Research Center for Particle Science and Technology, Institute of Frontier and Interdisciplinary Science, Shandong University, Qingdao, 266237, China; Research Center for Particle Science and Technology, Institute of Frontier and Interdisciplinary Science, Shandong University, Qingdao, 266237, China; School of Information Science and Engineering, Harbin Institute of Technology, Weihai, 264209, China; Research Center for Particle Science and Technology, Institute of Frontier and Interdisciplinary Science, Shandong University, Qingdao, 266237, China; School of Physics, Nanjing University, Nanjing, 210093, China; Research Center for Particle Science and Technology, Institute of Frontier and Interdisciplinary Science, Shandong University, Qingdao, 266237, China; School of Information Science and Engineering, Harbin Institute of Technology, Weihai, 264209, China; Institute of High Energy Physics, Chinese Academy of Sciences, 19B Yuquan Road, Shijingshan District, Beijing, 100049, China; School of Physics, Nankai University, Tianjin, 300071, China; Institute of High Energy Physics, Chinese Academy of Sciences, 19B Yuquan Road, Shijingshan District, Beijing, 100049, China'
$ws.Range("C43").Value = 'https://openalex.org/W4389617103'
$ws.Range("D43").Value = 'Simulation study of BESIII with stitched CMOS pixel detector using acts'
$ws.Range("E43").Value = '2023-12-01'
$ws.Range("F43").Value = 'Nuclear Science and Techniques'
$ws.Range("G43").Value = 'Springer Nature'
$ws.Range("H43").Value = 'https://doi.org/10.1007/s41365-023-01353-6'
$ws.Range("I43").Value = 'N/A'
$ws.Range("J43").Value = 'N/A'
$ws.Range("K43").Value = 'closed'
$ws.Range("L43").Value = 'en'
$ws.Range("M43").Value = '0'
$ws.Range("N43").Value = '2023'
$ws.Range("O43").Value = 'NA'
$ws.Range("P43").Value = 'https://doi.org/10.1007/s41365-023-01353-6'
$ws.Range("Q43").Value = 'article'

# Row 44
$ws.Range("A44").Value = 'Fan Wang, Moulun Luo'
$ws.Range("B44").Value = '; '
$ws.Range("C44").Value = 'https://openalex.org/W4390080872'
$ws.Range("D44").Value = 'Structural View of Human ABCs in Multidrug Resistance'
$ws.Range("E44").Value = '2023-12-21'
$ws.Range("F44").Value = 'N/A'
$ws.Range("G44").Value = 'N/A'
$ws.Range("H44").Value = 'https://doi.org/10.20944/preprints202312.1514.v1'
$ws.Range("I44").Value = 'N/A'
$ws.Range("J44").Value = 'submittedVersion'
$ws.Range("K44").Value = 'bronze'
$ws.Range("L44").Value = 'en'
$ws.Range("M44").Value = '0'
$ws.Range("N44").Value = '2023'
$ws.Range("O44").Value = 'NA'
$ws.Range("P44").Value = 'https://doi.org/10.20944/preprints202312.1514.v1'
$ws.Range("Q44").Value = 'article'

# Row 45
$ws.Range("A45").Value = 'Yuan Xiao Zhu, Chang‐Xin Shi, Laura A. Bruins, Patrick Jedlowski, Xuewei Wang, K. Martin Kortüm, Moulun Luo, Jonathan M. Ahmann, Esteban Braggio, A. Keith Stewart'
$ws.Range("B45").Value = '; ; ; ; ; ; ; ; ; '
$ws.Range("C45").Value = 'https://openalex.org/W4392724267'
$ws.Range("D45").Value = 'Data from Loss of &lt;i&gt;FAM46C&lt;/i&gt; Promotes Cell Survival in Myeloma'
$ws.Range("E45").Value = '2023-03-31'
$ws.Range("F45").Value = 'N/A'
$ws.Range("G45").Value = 'N/A'
$ws.Range("H45").Value = 'https://doi.org/10.1158/0008-5472.c.6509109'
$ws.Range("I45").Value = 'N/A'
$ws.Range("J45").Value = 'submittedVersion'
$ws.Range("K45").Value = 'closed'
$ws.Range("L45").Value = 'en'
$ws.Range("M45").Value = '0'
$ws.Range("N45").Value = '2023'
$ws.Range("O45").Value = 'NA'
$ws.Range("P45").Value = 'https://doi.org/10.1158/0008-5472.c.6509109'
$ws.Range("Q45").Value = 'article'
